# "Update Specification & Plannification"
#
# The task-planning table on sheet "Feuil1" is restructured: task numbers
# are renumbered, several task descriptions are reworded, some rows are
# dropped (e.g. "Tests", "Communication I2C" alone, "T90.10", ...), and a
# number of new sub-tasks are introduced (power/command/supply schematics,
# three-phase bridge management, UART/I2C-SPI communication, weekly
# supervisor meetings, ...). Rebuild the sheet from scratch with the new
# layout, values and formula, then re-apply the bold section-header /
# left-aligned sub-row formatting and cell merges.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Start clean: drop old merges & all cell content/formatting ---
$ws.Cells.UnMerge()
$ws.Cells.Clear()

# --- Table values (existing labels / numbers / estimates + totals formula) -
$ws.Range("A1").Value = "N°"
$ws.Range("B1").Value = "Description"
$ws.Range("D1").Value = "Estimation [h]"
$ws.Range("A2").Value = "T10"
$ws.Range("B2").Value = "Analyse des documents existants et recherches"
$ws.Range("A3").Value = "T10.10"
$ws.Range("C3").Value = "Analyse du schéma de la carte électronique de 2008"
$ws.Range("D3").Value = 5
$ws.Range("A4").Value = "T10.20"
$ws.Range("C4").Value = "Analyse du PCB de la carte électronique de 2008"
$ws.Range("D4").Value = 5
$ws.Range("A5").Value = "T10.30"
$ws.Range("C5").Value = "Analyse du code source de la carte électronique de 2008"
$ws.Range("D5").Value = 5
$ws.Range("A6").Value = "T10.40"
$ws.Range("C6").Value = "Recherches liées au projet"
$ws.Range("D6").Value = 15
$ws.Range("A7").Value = "T20"
$ws.Range("B7").Value = "Dessin du schéma électrique"
$ws.Range("D8").Value = 15
$ws.Range("D9").Value = 15
$ws.Range("D10").Value = 15
$ws.Range("A11").Value = "T30"
$ws.Range("B11").Value = "Dessin du (ou des) PCB(s)"
$ws.Range("D11").Value = 80
$ws.Range("A13").Value = "T40"
$ws.Range("B13").Value = "Commandes du matériel"
$ws.Range("A14").Value = "T40.10"
$ws.Range("C14").Value = "Faire la liste complète du matériel"
$ws.Range("D14").Value = 3
$ws.Range("A15").Value = "T40.20"
$ws.Range("C15").Value = "Rechercher le matériel sur internet"
$ws.Range("D15").Value = 20
$ws.Range("A16").Value = "T40.30"
$ws.Range("C16").Value = "Commander le matériel"
$ws.Range("D16").Value = 5
$ws.Range("A17").Value = "T50"
$ws.Range("B17").Value = "Programmation"
$ws.Range("A18").Value = "T50.10"
$ws.Range("D18").Value = 20
$ws.Range("A19").Value = "T50.20"
$ws.Range("C19").Value = "Gestion des capteurs à effet Hall digitaux"
$ws.Range("D19").Value = 10
$ws.Range("A20").Value = "T50.30"
$ws.Range("C20").Value = "Gestion des capteurs à effet Hall analogiques"
$ws.Range("D20").Value = 10
$ws.Range("A21").Value = "T50.40"
$ws.Range("C21").Value = "Gestion d'un codeur incrémental digital"
$ws.Range("D21").Value = 10
$ws.Range("A22").Value = "T50.50"
$ws.Range("C22").Value = "Gestion d'un capteur SIN/COS"
$ws.Range("D22").Value = 10
$ws.Range("A23").Value = "T50.60"
$ws.Range("C23").Value = "Gestion d'un codeur absolu digital"
$ws.Range("D23").Value = 10
$ws.Range("A24").Value = "T50.70"
$ws.Range("C24").Value = "Gestion d'un codeur absolu analogique"
$ws.Range("D24").Value = 10
$ws.Range("A25").Value = "T50.80"
$ws.Range("D25").Value = 10
$ws.Range("C26").Value = "Communication CAN"
$ws.Range("D26").Value = 25
$ws.Range("D27").Value = 10
$ws.Range("A28").Value = "T60"
$ws.Range("B28").Value = "Mesures des performances"
$ws.Range("D28").Value = 10
$ws.Range("A30").Value = "T70"
$ws.Range("B30").Value = "Rédaction de documents livrables"
$ws.Range("C31").Value = "Écriture du rapport"
$ws.Range("D31").Value = 30
$ws.Range("C32").Value = "Dessin des diagrammes UML"
$ws.Range("D32").Value = 10
$ws.Range("A33").Value = "T80"
$ws.Range("A34").Value = "T80.10"
$ws.Range("C34").Value = "Tâches administratives"
$ws.Range("D34").Value = 10
$ws.Range("A35").Value = "T80.20"
$ws.Range("D35").Value = 10
$ws.Range("A36").Value = "T90"
$ws.Range("B36").Value = "Imprévus (~10% = 42H)"
$ws.Range("D36").Value = 42
$ws.Range("D38").Formula = "=SUM(D2:D36)"

# --- Newly introduced rows/labels (written last so they land at the end of
# the shared-string table, mirroring how they were authored) -----------------
$ws.Range("C25").Value = "Communication I2C / SPI"
$ws.Range("A26").Value = "T50.90"
$ws.Range("C27").Value = "Communication UART"
$ws.Range("A27").Value = "T50.100"
$ws.Range("C18").Value = "Gestion du pont triphasé"
$ws.Range("A31").Value = "T70.10"
$ws.Range("A32").Value = "T70.20"
$ws.Range("C35").Value = "Séances ébdomadaires avec prof. Répondant (préparation)"
$ws.Range("B33").Value = "Gestion du projet"
$ws.Range("C8").Value = "Schéma de puissance"
$ws.Range("C9").Value = "Schéma de commande"
$ws.Range("C10").Value = "Schéma d'alimentations"

# --- Formatting -------------------------------------------------------------
# Row 1 (column headers): merged B:C, left aligned
$ws.Range("B1:C1").Merge()
$ws.Range("B1:C1").HorizontalAlignment = -4131

# Section-header rows: bold "N°" cell, merged+bold+left "Description" cell
$sectionRows = @(2, 7, 11, 13, 17, 28, 30, 33, 36)
foreach ($r in $sectionRows) {
    $ws.Range("A$r").Font.Bold = $true
    $bc = $ws.Range("B" + $r + ":C" + $r)
    $bc.Merge()
    $bc.Font.Bold = $true
    $bc.HorizontalAlignment = -4131
}

# Sub-rows with an indented (non-merged) left-aligned B:C pair
$leftAlignRows = @(8, 9, 10, 12, 18, 19, 29)
foreach ($r in $leftAlignRows) {
    $ws.Range("B" + $r + ":C" + $r).HorizontalAlignment = -4131
}

# --- View state ---------------------------------------------------------
$ws.Range("G12").Select()
try {
    $excel.ActiveWindow.ScrollRow = 8
    $excel.ActiveWindow.ScrollColumn = 1
} catch {}

